$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.628.91"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "1.794.88"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'227.12"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "'0.559"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("D8").Value = "'32.97"
$ws.Range("E8").Value = "  +3.89%  "
$ws.Range("D9").Value = "'0.298"
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "2.053.77"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.813.88"
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'11.09"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").Value = "'0.638"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").Value = "34.570.33"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").Value = "'4.30"
$ws.Range("E17").Value = "  +2.87%  "
$ws.Range("D18").Value = "'68.84"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").Value = "'248.41"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").Value = "0.0₃0801"
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("D21").Value = "'11.29"
$ws.Range("E21").Value = "  +3.01%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").Value = "'165.73"
$ws.Range("E25").Value = "  +2.40%  "
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D28").Value = "'0.117"
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'4.13"
$ws.Range("E30").Value = "  +13.55%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0526"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.24"
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'3.82"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("E34").Value = "  +2.86%  "
$ws.Range("D35").Value = "1.428.68"
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("E36").Value = "  +6.44%  "
$ws.Range("E37").Value = "  +2.93%  "
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("D40").Value = "'85.33"
$ws.Range("E40").Value = "  +6.52%  "
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("D43").Value = "'2.76"
$ws.Range("E43").Value = "  +3.01%  "
$ws.Range("D44").Value = "'13.73"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("E45").Value = "  +3.78%  "
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "1.953.64"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "'106.13"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("E51").Value = "  -6.83%  "
